$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Allow members to register" backlog item (row 2, Priority 0) is being
# re-prioritised: it moves down to sit immediately before the Priority-1
# group, and the items that used to follow it shift up to take its place.

# Capture the row-2 content before it is removed.
$taskB = $ws.Range("B2").Value()
$taskC = $ws.Range("C2").Value()
$taskD = $ws.Range("D2").Value()

# Remove row 2 entirely; rows 3-6 shift up to rows 2-5.
$ws.Rows(2).Delete()

# Insert a fresh row at the new position (row 6, just above the old row 7 /
# now-shifted Priority-1 group) and drop the relocated task there with its
# priority bumped up to 1.
$ws.Rows(6).Insert()
$ws.Range("A6").Value = 1.0
$ws.Range("B6").Value = $taskB
$ws.Range("C6").Value = $taskC
$ws.Range("D6").Value = $taskD
